# Generate Report for Handback
# Update the "last generated" timestamp strings on the Overview, zh-cn and
# de-de worksheets to reflect a fresh handback report run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G2)
$wsOverview.Range("G2").Value = "2016-08-31 09:20:56"

# zh-cn sheet: "Correspond Handoff Datetime" (H2) and
# "Correspond Handback DateTime" (K2)
$wsZhCn.Range("H2").Value = "2016-08-31 09:20:51"
$wsZhCn.Range("K2").Value = "2016-08-31 09:21:24"

# de-de sheet: "Correspond Handback DateTime" (K2)
$wsDeDe.Range("K2").Value = "2016-08-31 09:21:31"
